$wb = $excel.ActiveWorkbook

# hunk 0 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1476.6522  # H40 (was 1607.1428)
$ws.Cells.Item(40, 9).Value = 1556.1538  # I40 (was 1765)
$ws.Cells.Item(40, 10).Value = 1373.3  # J40 (was 1463.6364)
$ws.Cells.Item(40, 11).Value = 1556.1538  # K40 (was 1765)
$ws.Cells.Item(40, 12).Value = 1373.3  # L40 (was 1463.6364)
$ws.Cells.Item(40, 13).Value = -1381.1538  # M40 (was -1590)
$ws.Cells.Item(40, 14).Value = -1723.3  # N40 (was -1813.6364)

# hunk 1 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 7752.381  # H131 (was 6111.069)
$ws.Cells.Item(131, 9).Value = 1050  # I131 (was 946.7778)
$ws.Cells.Item(131, 10).Value = 9329.412  # J131 (was 8435)
$ws.Cells.Item(131, 11).Value = 3150  # K131 (was 2840.3334)
$ws.Cells.Item(131, 12).Value = 27988.236  # L131 (was 25305)
$ws.Cells.Item(131, 13).Value = 1890  # M131 (was 2199.6666)
$ws.Cells.Item(131, 14).Value = -38068.236  # N131 (was -35385)

# hunk 2 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3849601.5  # H132 (was 4314157.5)
$ws.Cells.Item(132, 9).Value = 3346.7166  # I132 (was 3732.7925)
$ws.Cells.Item(132, 11).Value = 10040.1498  # K132 (was 11198.3775)
$ws.Cells.Item(132, 13).Value = -7510.149800000001  # M132 (was -8668.377500000001)

# hunk 3 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 1224.6875  # H141 (was 1107.7)
$ws.Cells.Item(141, 9).Value = 1167.4193  # I141 (was 1001.73914)
$ws.Cells.Item(141, 10).Value = 3000  # J141 (was 2326.25)
$ws.Cells.Item(141, 11).Value = 3502.2579  # K141 (was 3005.21742)
$ws.Cells.Item(141, 12).Value = 9000  # L141 (was 6978.75)
$ws.Cells.Item(141, 13).Value = 1677.7421  # M141 (was 2174.78258)
$ws.Cells.Item(141, 14).Value = -19360  # N141 (was -17338.75)

# hunk 4 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5558.54  # H32 (was 8429.885)
$ws.Cells.Item(32, 9).Value = 5254.6377  # I32 (was 8807.308999999999)
$ws.Cells.Item(32, 10).Value = 6774.15  # J32 (was 7079.1055)
$ws.Cells.Item(32, 11).Value = 5254.6377  # K32 (was 8807.308999999999)
$ws.Cells.Item(32, 12).Value = 6774.15  # L32 (was 7079.1055)
$ws.Cells.Item(32, 13).Value = -4967.6377  # M32 (was -8520.308999999999)
$ws.Cells.Item(32, 14).Value = -7348.15  # N32 (was -7653.1055)

# hunk 5 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 8621940  # H61 (was 12501597)
$ws.Cells.Item(61, 9).Value = 9435254  # I61 (was 13159531)
$ws.Cells.Item(61, 10).Value = 805.6  # J61 (was 864)
$ws.Cells.Item(61, 11).Value = 9435254  # K61 (was 13159531)
$ws.Cells.Item(61, 12).Value = 805.6  # L61 (was 864)
$ws.Cells.Item(61, 13).Value = -9435042  # M61 (was -13159319)
$ws.Cells.Item(61, 14).Value = -1229.6  # N61 (was -1288)

# hunk 6 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 9616853  # H74 (was 7043353.5)
$ws.Cells.Item(74, 9).Value = 11906151  # I74 (was 9435074)
$ws.Cells.Item(74, 10).Value = 1799.4  # J74 (was 1065.9445)
$ws.Cells.Item(74, 11).Value = 11906151  # K74 (was 9435074)
$ws.Cells.Item(74, 12).Value = 1799.4  # L74 (was 1065.9445)
$ws.Cells.Item(74, 13).Value = -11905277  # M74 (was -9434200)
$ws.Cells.Item(74, 14).Value = -3547.4  # N74 (was -2813.9445)

# hunk 7 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 9616853  # H77 (was 7043353.5)
$ws.Cells.Item(77, 9).Value = 11906151  # I77 (was 9435074)
$ws.Cells.Item(77, 10).Value = 1799.4  # J77 (was 1065.9445)
$ws.Cells.Item(77, 11).Value = 59530755  # K77 (was 47175370)
$ws.Cells.Item(77, 12).Value = 8997  # L77 (was 5329.7225)
$ws.Cells.Item(77, 13).Value = -59526387  # M77 (was -47171002)
$ws.Cells.Item(77, 14).Value = -17733  # N77 (was -14065.7225)

# hunk 8 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4809366  # H132 (was 5320960)
$ws.Cells.Item(132, 9).Value = 5320638.5  # I132 (was 5954002)
$ws.Cells.Item(132, 11).Value = 15961915.5  # K132 (was 17862006)
$ws.Cells.Item(132, 13).Value = -15959385.5  # M132 (was -17859476)

# hunk 9 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 8621940  # H136 (was 12501597)
$ws.Cells.Item(136, 9).Value = 9435254  # I136 (was 13159531)
$ws.Cells.Item(136, 10).Value = 805.6  # J136 (was 864)
$ws.Cells.Item(136, 11).Value = 28305762  # K136 (was 39478593)
$ws.Cells.Item(136, 12).Value = 2416.8  # L136 (was 2592)
$ws.Cells.Item(136, 13).Value = -28303212  # M136 (was -39476043)
$ws.Cells.Item(136, 14).Value = -7516.8  # N136 (was -7692)

# hunk 10 -> sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 0  # H2 (was 67114.28999999999)
$ws.Cells.Item(2, 10).Value = 0  # J2 (was 67114.28999999999)
$ws.Cells.Item(2, 12).ClearContents()  # L2 (was 67114.28999999999)
$ws.Cells.Item(2, 14).Value = 0  # N2 (was -67340.28999999999)

# hunk 11 -> sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2105.016  # H134 (was 2454.7886)
$ws.Cells.Item(134, 9).Value = 1368.721  # I134 (was 1560.8918)
$ws.Cells.Item(134, 10).Value = 3771.3684  # J134 (was 4659.7334)
$ws.Cells.Item(134, 11).Value = 4106.163  # K134 (was 4682.6754)
$ws.Cells.Item(134, 12).Value = 11314.1052  # L134 (was 13979.2002)
$ws.Cells.Item(134, 13).Value = -1571.163  # M134 (was -2147.6754)
$ws.Cells.Item(134, 14).Value = -16384.1052  # N134 (was -19049.2002)

# hunk 12 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5130060.5  # H31 (was 6412629.5)
$ws.Cells.Item(31, 9).Value = 1542.3726  # I31 (was 2089)
$ws.Cells.Item(31, 10).Value = 23812520  # J31 (was 25644252)
$ws.Cells.Item(31, 11).Value = 1542.3726  # K31 (was 2089)
$ws.Cells.Item(31, 12).Value = 23812520  # L31 (was 25644252)
$ws.Cells.Item(31, 13).Value = -1247.3726  # M31 (was -1794)
$ws.Cells.Item(31, 14).Value = -23813110  # N31 (was -25644842)

# hunk 13 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5130060.5  # H34 (was 6412629.5)
$ws.Cells.Item(34, 9).Value = 1542.3726  # I34 (was 2089)
$ws.Cells.Item(34, 10).Value = 23812520  # J34 (was 25644252)
$ws.Cells.Item(34, 11).Value = 1542.3726  # K34 (was 2089)
$ws.Cells.Item(34, 12).Value = 23812520  # L34 (was 25644252)
$ws.Cells.Item(34, 13).Value = -1340.3726  # M34 (was -1887)
$ws.Cells.Item(34, 14).Value = -23812924  # N34 (was -25644656)

# hunk 14 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(81, 8).Value = 48996.8  # H81 (was 36496)
$ws.Cells.Item(81, 10).Value = 48996.8  # J81 (was 36496)
$ws.Cells.Item(81, 12).Value = 48996.8  # L81 (was 36496)
$ws.Cells.Item(81, 14).Value = -50992.8  # N81 (was -38492)

# hunk 15 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(84, 8).Value = 48996.8  # H84 (was 36496)
$ws.Cells.Item(84, 10).Value = 48996.8  # J84 (was 36496)
$ws.Cells.Item(84, 12).Value = 146990.4  # L84 (was 109488)
$ws.Cells.Item(84, 14).Value = -156974.4  # N84 (was -119472)

# hunk 16 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1477.3334  # H99 (was 1332.6666)
$ws.Cells.Item(99, 9).Value = 1532.5  # I99 (was 1178.2106)
$ws.Cells.Item(99, 10).Value = 1367  # J99 (was 2800)
$ws.Cells.Item(99, 11).Value = 1532.5  # K99 (was 1178.2106)
$ws.Cells.Item(99, 12).Value = 1367  # L99 (was 2800)
$ws.Cells.Item(99, 13).Value = -34.5  # M99 (was 319.7893999999999)
$ws.Cells.Item(99, 14).Value = -4363  # N99 (was -5796)

# hunk 17 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1477.3334  # H126 (was 1332.6666)
$ws.Cells.Item(126, 9).Value = 1532.5  # I126 (was 1178.2106)
$ws.Cells.Item(126, 10).Value = 1367  # J126 (was 2800)
$ws.Cells.Item(126, 11).Value = 4597.5  # K126 (was 3534.6318)
$ws.Cells.Item(126, 12).Value = 4101  # L126 (was 8400)
$ws.Cells.Item(126, 13).Value = -2127.5  # M126 (was -1064.6318)
$ws.Cells.Item(126, 14).Value = -9041  # N126 (was -13340)

# hunk 18 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 776.1905  # H5 (was 902.94116)
$ws.Cells.Item(5, 9).Value = 340.75  # I5 (was 392.375)
$ws.Cells.Item(5, 11).Value = 1022.25  # K5 (was 1177.125)
$ws.Cells.Item(5, 13).Value = -910.25  # M5 (was -1065.125)

# hunk 19 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value = 9675.666999999999  # H110 (was 8506.5)
$ws.Cells.Item(110, 9).Value = 4513.5  # I110 (was 4008.6667)
$ws.Cells.Item(110, 10).Value = 20000  # J110 (was 22000)
$ws.Cells.Item(110, 11).Value = 13540.5  # K110 (was 12026.0001)
$ws.Cells.Item(110, 12).Value = 60000  # L110 (was 66000)
$ws.Cells.Item(110, 13).Value = -9450.5  # M110 (was -7936.000100000001)
$ws.Cells.Item(110, 14).Value = -68180  # N110 (was -74180)

# hunk 20 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 4546.4688  # H112 (was 4849.512)
$ws.Cells.Item(112, 9).Value = 6702.6665  # I112 (was 5776)
$ws.Cells.Item(112, 10).Value = 4323.4136  # J112 (was 4720.8335)
$ws.Cells.Item(112, 11).Value = 20107.9995  # K112 (was 17328)
$ws.Cells.Item(112, 12).Value = 12970.2408  # L112 (was 14162.5005)
$ws.Cells.Item(112, 13).Value = -18999.9995  # M112 (was -16220)
$ws.Cells.Item(112, 14).Value = -15186.2408  # N112 (was -16378.5005)

# hunk 21 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 713.3291  # H113 (was 699.35364)
$ws.Cells.Item(113, 10).Value = 868.92  # J113 (was 838.49054)
$ws.Cells.Item(113, 12).Value = 2606.76  # L113 (was 2515.47162)
$ws.Cells.Item(113, 14).Value = -6946.76  # N113 (was -6855.47162)

# hunk 22 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(116, 8).Value = 4178.2  # H116 (was 615.6)
$ws.Cells.Item(116, 9).Value = 1029  # I116 (was 519.5)
$ws.Cells.Item(116, 10).Value = 4965.5  # J116 (was 1000)
$ws.Cells.Item(116, 11).Value = 3087  # K116 (was 1558.5)
$ws.Cells.Item(116, 12).Value = 14896.5  # L116 (was 3000)
$ws.Cells.Item(116, 13).Value = 355  # M116 (was 1883.5)
$ws.Cells.Item(116, 14).Value = -21780.5  # N116 (was -9884)

# hunk 23 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(118, 8).Value = 1258.2727  # H118 (was 1442.6)
$ws.Cells.Item(118, 9).Value = 1500  # I118 (was 985.6667)
$ws.Cells.Item(118, 10).Value = 1246.762  # J118 (was 1556.8334)
$ws.Cells.Item(118, 11).Value = 4500  # K118 (was 2957.0001)
$ws.Cells.Item(118, 12).Value = 3740.286  # L118 (was 4670.5002)
$ws.Cells.Item(118, 13).Value = -3257  # M118 (was -1714.0001)
$ws.Cells.Item(118, 14).Value = -6226.286  # N118 (was -7156.5002)

# hunk 24 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1534.8182  # H122 (was 1711.4445)
$ws.Cells.Item(122, 10).Value = 1017  # J122 (was 1201.6666)
$ws.Cells.Item(122, 12).Value = 9153  # L122 (was 10814.9994)
$ws.Cells.Item(122, 14).Value = -14053  # N122 (was -15714.9994)

# hunk 25 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 776.1905  # H135 (was 902.94116)
$ws.Cells.Item(135, 9).Value = 340.75  # I135 (was 392.375)
$ws.Cells.Item(135, 11).Value = 3066.75  # K135 (was 3531.375)
$ws.Cells.Item(135, 13).Value = -531.75  # M135 (was -996.375)

# hunk 26 -> sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 59454.4  # H104 (was 60503)
$ws.Cells.Item(104, 10).Value = 59454.4  # J104 (was 60503)
$ws.Cells.Item(104, 12).Value = 59454.4  # L104 (was 60503)
$ws.Cells.Item(104, 14).Value = -66442.39999999999  # N104 (was -67491)

# hunk 27 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1354.84  # H22 (was 1264.4073)
$ws.Cells.Item(22, 9).Value = 540  # I22 (was 416.25)
$ws.Cells.Item(22, 10).Value = 1558.55  # J22 (was 1621.5264)
$ws.Cells.Item(22, 11).Value = 540  # K22 (was 416.25)
$ws.Cells.Item(22, 12).Value = 1558.55  # L22 (was 1621.5264)
$ws.Cells.Item(22, 13).Value = -245  # M22 (was -121.25)
$ws.Cells.Item(22, 14).Value = -2148.55  # N22 (was -2211.5264)

# hunk 28 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 1354.84  # H27 (was 1264.4073)
$ws.Cells.Item(27, 9).Value = 540  # I27 (was 416.25)
$ws.Cells.Item(27, 10).Value = 1558.55  # J27 (was 1621.5264)
$ws.Cells.Item(27, 11).Value = 540  # K27 (was 416.25)
$ws.Cells.Item(27, 12).Value = 1558.55  # L27 (was 1621.5264)
$ws.Cells.Item(27, 13).Value = -433  # M27 (was -309.25)
$ws.Cells.Item(27, 14).Value = -1772.55  # N27 (was -1835.5264)

# hunk 29 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 877.1177  # H46 (was 1073.9)
$ws.Cells.Item(46, 9).Value = 570.2  # I46 (was 790)
$ws.Cells.Item(46, 10).Value = 1005  # J46 (was 1144.875)
$ws.Cells.Item(46, 11).Value = 570.2  # K46 (was 790)
$ws.Cells.Item(46, 12).Value = 1005  # L46 (was 1144.875)
$ws.Cells.Item(46, 13).Value = -382.2  # M46 (was -602)
$ws.Cells.Item(46, 14).Value = -1381  # N46 (was -1520.875)

# hunk 30 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(117, 8).Value = 50392  # H117 (was 0)
$ws.Cells.Item(117, 10).Value = 50392  # J117 (was 0)
$ws.Cells.Item(117, 12).Value = 50392  # L117 (was 0)
$ws.Cells.Item(117, 14).Value = -59570  # N117 (was None)

# hunk 31 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 10878945  # H132 (was 15637005)
$ws.Cells.Item(132, 9).Value = 6567.92  # I132 (was 9342.714)
$ws.Cells.Item(132, 10).Value = 23822250  # J132 (was 27791854)
$ws.Cells.Item(132, 11).Value = 19703.76  # K132 (was 28028.142)
$ws.Cells.Item(132, 12).Value = 71466750  # L132 (was 83375562)
$ws.Cells.Item(132, 13).Value = -17173.76  # M132 (was -25498.142)
$ws.Cells.Item(132, 14).Value = -71471810  # N132 (was -83380622)

# hunk 32 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 55550  # H139 (was 55716.668)
$ws.Cells.Item(139, 10).Value = 55550  # J139 (was 55716.668)
$ws.Cells.Item(139, 12).Value = 55550  # L139 (was 55716.668)
$ws.Cells.Item(139, 14).Value = -65830  # N139 (was -65996.66800000001)

# hunk 33 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 930.5682  # H136 (was 991.2195)
$ws.Cells.Item(136, 9).Value = 722.8421  # I136 (was 776.25714)
$ws.Cells.Item(136, 10).Value = 2246.1667  # J136 (was 2245.1667)
$ws.Cells.Item(136, 11).Value = 2168.5263  # K136 (was 2328.77142)
$ws.Cells.Item(136, 12).Value = 6738.500100000001  # L136 (was 6735.500100000001)
$ws.Cells.Item(136, 13).Value = 381.4737  # M136 (was 221.22858)
$ws.Cells.Item(136, 14).Value = -11838.5001  # N136 (was -11835.5001)
